$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add new row 8 of trade data
$ws.Range("A8").Value = 8256.42
$ws.Range("B8").Value = 8145.64
$ws.Range("C8").Value = 19.170000000000002
$ws.Range("D8").Value = 19.43
$ws.Range("E8").Value = $false
$ws.Range("F8").Value = 1.36
$ws.Range("G8").Value = 42609.488749999997
$ws.Range("H8").Value = $true

# Match the date/time number format already used by column G (style index 1)
$ws.Range("G8").NumberFormat = "m/d/yy h:mm"
